# === Reorganize workbook: add 'Player Info' sheet before 'ODI Batting',
# add 'ODI Batting Extra' sheet after 'ODI Bowling', rename link columns
# to match-code columns, and tidy up blank inning-number cells. ===

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# --- New first sheet: Player Info ---
$wsPlayer = $wb.Worksheets.Add($wsBatting)
$wsPlayer.Name = "Player Info"

$playerHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $playerHeaders.Length; $i++) {
    $cell = $wsPlayer.Cells.Item(1, $i + 1)
    $cell.Value = $playerHeaders[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$wsPlayer.Cells.Item(2, 1).NumberFormat = "@"
$wsPlayer.Cells.Item(2, 1).Value = "4652"
$wsPlayer.Cells.Item(2, 2).Value = "Kuldeep Yadav"
$wsPlayer.Cells.Item(2, 3).Value = "Left Handed"
$wsPlayer.Cells.Item(2, 4).Value = "Left Arm Wrist Spin (Chinaman)"

# --- New last sheet: ODI Batting Extra ---
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Length; $i++) {
    $cell = $wsExtra.Cells.Item(1, $i + 1)
    $cell.Value = $extraHeaders[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$extraData = @(
    ("4454", 10, $null, $null, $null, "NO")
    ,("4456", $null, $null, $null, $null, "NO")
    ,("4480", 10, $null, $null, $null, "NO")
    ,("4482", 10, $null, $null, $null, "NO")
    ,("4536", 9, "0", "0", "1.89%", "NO")
    ,("4637", 9, $null, $null, $null, "NO")
    ,("4640", 9, $null, $null, $null, "NO")
    ,("4643", $null, $null, $null, $null, "NO")
    ,("4656", $null, $null, $null, $null, "NO")
    ,("4657", 9, $null, $null, $null, "NO")
    ,("4658", 9, $null, $null, $null, "YES")
    ,("4685", $null, $null, $null, $null, "NO")
    ,("4689", 8, "2", "0", "4.57%", "YES")
    ,("4691", $null, $null, $null, $null, "NO")
    ,("4692", $null, $null, $null, $null, "NO")
    ,("4695", 9, $null, $null, $null, "NO")
    ,("4697", 9, "0", "0", "0.78%", "NO")
    ,("4725", 9, $null, $null, $null, "NO")
    ,("4728", 9, "0", "0", "3.42%", "NO")
    ,("4732", 9, "0", "0", "2.42%", "NO")
)

$r = 2
foreach ($row in $extraData) {
    $wsExtra.Cells.Item($r, 1).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $wsExtra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $wsExtra.Cells.Item($r, 3).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $wsExtra.Cells.Item($r, 4).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($null -ne $row[4]) {
        $wsExtra.Cells.Item($r, 5).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 5).Value = $row[4]
    }
    $wsExtra.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# --- ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE, replace URL with code ---
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingCodes = @{
    2 = "4051"
    3 = "4052"
    4 = "4053"
    5 = "4056"
    6 = "4057"
    7 = "4064"
    8 = "4065"
    9 = "4067"
    10 = "4069"
    11 = "4071"
    12 = "4076"
    13 = "4085"
    14 = "4096"
    15 = "4099"
    16 = "4126"
    17 = "4127"
    18 = "4128"
    19 = "4130"
    20 = "4133"
    21 = "4135"
    22 = "4171"
    23 = "4173"
    24 = "4175"
    25 = "4196"
    26 = "4197"
    27 = "4199"
    28 = "4201"
    29 = "4203"
    30 = "4205"
    31 = "4216"
    32 = "4219"
    33 = "4220"
    34 = "4221"
    35 = "4234"
    36 = "4235"
    37 = "4239"
    38 = "4242"
    39 = "4245"
    40 = "4248"
    41 = "4258"
    42 = "4263"
    43 = "4266"
    44 = "4268"
    45 = "4270"
    46 = "4310"
    47 = "4316"
    48 = "4324"
    49 = "4332"
    50 = "4338"
    51 = "4342"
    52 = "4350"
    53 = "4359"
    54 = "4360"
    55 = "4385"
    56 = "4387"
    57 = "4388"
    58 = "4398"
    59 = "4399"
    60 = "4400"
    61 = "4402"
    62 = "4437"
    63 = "4454"
    64 = "4456"
    65 = "4480"
    66 = "4482"
    67 = "4536"
    68 = "4637"
    69 = "4640"
    70 = "4643"
    71 = "4656"
    72 = "4657"
    73 = "4658"
    74 = "4685"
    75 = "4689"
    76 = "4691"
    77 = "4692"
    78 = "4695"
    79 = "4697"
    80 = "4725"
    81 = "4728"
    82 = "4732"
}
foreach ($r in $battingCodes.Keys) {
    $wsBatting.Cells.Item($r, 4).NumberFormat = "@"
    $wsBatting.Cells.Item($r, 4).Value = $battingCodes[$r]
}

# --- ODI Batting sheet: clear blank INNING_NUMBER cells ---
$battingBlankInningRows = @(2, 3, 4, 6, 7, 8, 11, 12, 15, 16, 17, 18, 21, 22, 24, 26, 27, 28, 31, 33, 34, 36, 37, 38, 39, 41, 46, 47, 48, 51, 52, 53, 54, 55, 56, 57, 59, 60, 61, 62, 63, 64, 65, 66, 68, 69, 72, 73, 76, 78, 80)
foreach ($r in $battingBlankInningRows) {
    $wsBatting.Cells.Item($r, 2).ClearContents()
}

# --- ODI Bowling sheet: rename MATCH_CARD_LINK -> MATCH_CODE, replace URL with code ---
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4052"
    3 = "4053"
    4 = "4056"
    5 = "4057"
    6 = "4064"
    7 = "4065"
    8 = "4067"
    9 = "4069"
    10 = "4071"
    11 = "4076"
    12 = "4085"
    13 = "4099"
    14 = "4126"
    15 = "4127"
    16 = "4128"
    17 = "4130"
    18 = "4133"
    19 = "4135"
    20 = "4171"
    21 = "4173"
    22 = "4175"
    23 = "4196"
    24 = "4197"
    25 = "4199"
    26 = "4201"
    27 = "4203"
    28 = "4205"
    29 = "4216"
    30 = "4219"
    31 = "4220"
    32 = "4221"
    33 = "4234"
    34 = "4235"
    35 = "4239"
    36 = "4242"
    37 = "4245"
    38 = "4248"
    39 = "4258"
    40 = "4263"
    41 = "4266"
    42 = "4268"
    43 = "4270"
    44 = "4310"
    45 = "4316"
    46 = "4324"
    47 = "4332"
    48 = "4338"
    49 = "4342"
    50 = "4350"
    51 = "4359"
    52 = "4360"
    53 = "4385"
    54 = "4387"
    55 = "4388"
    56 = "4398"
    57 = "4399"
    58 = "4400"
    59 = "4402"
    60 = "4437"
    61 = "4454"
    62 = "4456"
    63 = "4480"
    64 = "4482"
    65 = "4536"
    66 = "4637"
    67 = "4640"
    68 = "4643"
    69 = "4656"
    70 = "4657"
    71 = "4658"
    72 = "4685"
    73 = "4689"
    74 = "4691"
    75 = "4692"
    76 = "4695"
    77 = "4697"
    78 = "4725"
    79 = "4728"
    80 = "4732"
}
foreach ($r in $bowlingCodes.Keys) {
    $wsBowling.Cells.Item($r, 2).NumberFormat = "@"
    $wsBowling.Cells.Item($r, 2).Value = $bowlingCodes[$r]
}

# --- Restore active sheet / selection to the new first sheet ---
$wsPlayer.Activate()
$wsPlayer.Range("A1").Select()

Write-Output "done"
